# Generate Report for Handback
#
# The workbook tracks localization handback status for source files.
# This edit:
#   1) Refreshes the existing "855bebd0-e14f-4b76-afd8-bfd13c3e8764.md"
#      entry's generated UUID (now "ab6cf907-ef68-4f0f-9ae1-a47026b937b8")
#      and its timestamps / derived xliff file names.
#   2) Appends a brand-new row for a second file,
#      "d2f17490-2dad-4f74-9f05-6ee478d42efb.md", across all three sheets
#      (Overview, zh-cn, de-de), including hyperlinks and table growth.

$wb = $excel.ActiveWorkbook

$oldGuid = "855bebd0-e14f-4b76-afd8-bfd13c3e8764"
$guid1   = "ab6cf907-ef68-4f0f-9ae1-a47026b937b8"
$guid2   = "d2f17490-2dad-4f74-9f05-6ee478d42efb"

$hash1 = "83c057c332cad23ee05fca9d24080b06bc355d72"
$hash2 = "c94ed268f6c6903dd1d3c94885ba664e6db4a140"

$srcCommit  = "c8a8df3dba688edf6ac190d4667facb86d47ce21"
$zhCommit   = "ed60780f3b775439bb6ad82dd8b97cc618595a01"
$deCommit   = "62db6d4653e46fba0e4f23665acf58f899c0c9b1"

function SrcUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$name"
}
function ZhUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhCommit/e2e/$name"
}
function DeUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$deCommit/e2e/$name"
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (A1:G.. ; table "Overview")
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item("Overview")

# --- update existing row 2 (guid1, was oldGuid) ---
$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Range("A2").Value = "$guid1.md"
$wsOv.Range("B2").Value = "e2e\$guid1.md"
$wsOv.Range("G2").Value = "2016-08-26 04:58:33"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), (SrcUrl "$guid1.md"), "", "", "e2e\$guid1.md")

# --- append new row 3 (guid2) ---
$loOv.ListRows.Add() | Out-Null
$wsOv.Range("A3").Value = "$guid2.md"
$wsOv.Range("B3").Value = "e2e\$guid2.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOv.Range("G3").Value = "2016-08-26 04:58:33"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), (SrcUrl "$guid2.md"), "", "", "e2e\$guid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (A1:P.. ; table "zh-cn")
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")

# --- update existing row 2 (guid1) ---
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("I2").Hyperlinks.Delete()

$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Range("G2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-26 04:58:28"
$wsZh.Range("I2").Value = "$guid1.md"
$wsZh.Range("J2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 04:58:45"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), (ZhUrl "$guid1.md"), "", "", "$guid1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), (ZhUrl "$guid1.md"), "", "", "$guid1.md")

# --- append new row 3 (guid2) ---
$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A3").Value = "$guid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = "2016-08-26 04:58:28"
$wsZh.Range("I3").Value = "$guid2.md"
$wsZh.Range("J3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "2016-08-26 04:58:45"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), (ZhUrl "$guid2.md"), "", "", "$guid2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), (ZhUrl "$guid2.md"), "", "", "$guid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (A1:P.. ; table "de-de")
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")

# --- update existing row 2 (guid1) ---
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("I2").Hyperlinks.Delete()

$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Range("G2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-26 04:58:33"
$wsDe.Range("I2").Value = "$guid1.md"
$wsDe.Range("J2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 04:58:52"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), (DeUrl "$guid1.md"), "", "", "$guid1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), (DeUrl "$guid1.md"), "", "", "$guid1.md")

# --- append new row 3 (guid2) ---
$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A3").Value = "$guid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = "2016-08-26 04:58:33"
$wsDe.Range("I3").Value = "$guid2.md"
$wsDe.Range("J3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "2016-08-26 04:58:52"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), (DeUrl "$guid2.md"), "", "", "$guid2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), (DeUrl "$guid2.md"), "", "", "$guid2.md")
